$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 (Enemy2): LIFE value changes from ">0" to ">10"
$ws.Range("B3").Value = ">10"

# Add new row 4: Enemy1, >10, <=100, TRUE
$ws.Range("A4").Value = "Enemy1"
$ws.Range("B4").Value = ">10"
$ws.Range("C4").Value = "<=100"
$ws.Range("D4").Value = $true

# Add new row 5: Enemy2, >100, >8, TRUE
$ws.Range("A5").Value = "Enemy2"
$ws.Range("B5").Value = ">100"
$ws.Range("C5").Value = ">8"
$ws.Range("D5").Value = $true

# Copy formatting from row 3 (existing styled row) to the new rows 4 and 5
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Range("A5:D5").PasteSpecial(-4122)
